$d = $word.ActiveDocument

# Ordered (old, new) pairs matching the document's top-to-bottom cell order.
# Doing the replacements in this order guarantees each "old" text is unique
# and unambiguous at the moment it is searched for (one pair's "new" text
# happens to equal another pair's "old" text, so order matters).
$pairs = @(
    @("87÷3=29, 0", "74÷4=18, 2"),
    @("12÷8=1, 4", "18÷5=3, 3"),
    @("97÷2=48, 1", "47÷3=15, 2"),
    @("72÷9=8, 0", "94÷2=47, 0"),
    @("38÷3=12, 2", "65÷9=7, 2"),
    @("16÷6=2, 4", "96÷6=16, 0"),
    @("62÷5=12, 2", "97÷7=13, 6"),
    @("60÷5=12, 0", "96÷2=48, 0"),
    @("71÷4=17, 3", "83÷8=10, 3"),
    @("12÷5=2, 2", "50÷5=10, 0"),
    @("27÷2=13, 1", "76÷6=12, 4"),
    @("89÷5=17, 4", "61÷6=10, 1"),
    @("51÷9=5, 6", "70÷2=35, 0"),
    @("75÷6=12, 3", "35÷7=5, 0"),
    @("24÷5=4, 4", "57÷5=11, 2"),
    @("93÷4=23, 1", "87÷3=29, 0"),
    @("32÷7=4, 4", "48÷6=8, 0"),
    @("58÷7=8, 2", "32÷8=4, 0"),
    @("52÷5=10, 2", "64÷7=9, 1"),
    @("70÷8=8, 6", "77÷2=38, 1"),
    @("69÷9=7, 6", "69÷5=13, 4"),
    @("17÷3=5, 2", "26÷3=8, 2"),
    @("30÷2=15, 0", "66÷6=11, 0"),
    @("34÷9=3, 7", "70÷4=17, 2"),
    @("98÷5=19, 3", "91÷4=22, 3")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
